$d = $word.ActiveDocument

# --- Locate the paragraph that currently holds the bare "[S19]" text ---
# (it sits right after the "[S18] ... le 16-12-2019." paragraph, which
#  currently carries the _GoBack bookmark).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r") -eq "[S19]") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the [S19] placeholder paragraph"
}

# --- Detach the _GoBack bookmark from its current (S18) location; it will
#     be re-inserted at the end of the new S19 paragraph below. ---
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Build the replacement paragraph (full [S19] reference entry) and
#     insert it right before the old placeholder paragraph. Word merges
#     a trailing empty <w:p/> back into what remains of the old
#     paragraph, so its own paragraph mark/properties are untouched and
#     the new paragraph keeps its own <w:pPr>. ---
$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
        '<w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr></w:pPr>' +
        '<w:r><w:t>[S</w:t></w:r>' +
        '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>19</w:t></w:r>' +
        '<w:r><w:t>]</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:fldChar w:fldCharType="begin"/></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:instrText xml:space="preserve"> HYPERLINK "https://fr.wikipedia.org/wiki/Client%E2%80%93serveur," </w:instrText></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:fldChar w:fldCharType="separate"/></w:r>' +
        '<w:r><w:rPr><w:rStyle w:val="7"/><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:t>https://fr.wikipedia.org/wiki/Client%E2%80%93serveur,</w:t></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:fldChar w:fldCharType="end"/></w:r>' +
        '<w:r><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="fr"/></w:rPr><w:t xml:space="preserve"> le 19-12-2019.</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '<w:p/>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml)
